$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary / header values -------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 928000
# Cant. Trabajadores (now a single worker)
$ws.Range("C13").Value = 1
# Cant. Periodos (now 20 periods)
$ws.Range("F13").Value = 20

# --- Detail rows --------------------------------------------------------------
# The table used to list 2 workers: "33101764 / DIANA KARINA ARIAS RIOS" with
# periods 2502..2507 (rows 16-21) and "1083032206 / MARIA ALEJANDRA ROA
# MEDINA" with periods 2401..2501 (rows 22-40). Only MARIA remains, covering
# periods 2401..2412 and 2501..2508 (20 rows total). Remove DIANA's 5 surplus
# rows first (row 21, her last row, already shares the regular row style, so
# deleting the first 5 of her 6 rows shifts everything up cleanly and leaves
# the special "last row" bottom-border formatting on the new final row, 35).
$ws.Range("16:20").Delete()

$periods = @("2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "1083032206"
    $ws.Cells.Item($r, 4).Value = "MARIA ALEJANDRA ROA MEDINA"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 46400
    $ws.Cells.Item($r, 7).Value = 1160000
}
